# Apply "repull data, push all data, mean calculation" changes:
# Update column F (dSF) values for rows where a newly computed mean-based
# value replaces the previous placeholder 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = -2
    4  = 3
    5  = -4
    6  = 7
    7  = -3
    9  = -5
    11 = 1
    12 = 1
    13 = -1
    14 = -1
    15 = 3
    16 = -2
    17 = -3
    18 = 4
    19 = -5
    20 = -3
    22 = -1
    25 = 3
    26 = 1
    27 = -1
    28 = -2
    29 = 4
    30 = -1
    31 = -5
    32 = -4
    33 = -2
    34 = -4
    36 = -1
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
